# The sheet tracked a table of anomaly-discord benchmark runs. This edit
# replaces the stale "nab-data ambient_temperature" / NAB row with a run
# against the UCR anomaly dataset (per commit: "test on UCR anomaly
# dataset but got score 0").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dataset name / length / #discords / discord position / discord length
$ws.Range("B2").Value = "227_UCR_Anomaly_mit14134longtermecg_11231_29000_29100.txt"
$ws.Range("C2").Value = 59302
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 29000
$ws.Range("G2").Value = 102

# Widen the dataset-name column to fit the new, longer file name, and give
# the newly-visible "Dataset length" column an explicit width too.
$ws.Columns("B").ColumnWidth = 58.5866666666667
$ws.Columns("C").ColumnWidth = 23.5266666666667

# Touch the last couple of rows (mirrors the row-height bookkeeping seen
# at the bottom of the sheet once it was scrolled/selected all the way
# down) and move the active selection to A2.
$ws.Rows(1048575).RowHeight = 12.8
$ws.Rows(1048576).RowHeight = 12.8
[void]$ws.Range("A2").Select()
